$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source row to copy the date-column (A) cell style/format from.
$styleSource = $ws.Range("A906")

$data = @(
    @(907, 45461, 0.13778, 0.13782, 0.13778, 0.13782, 0),
    @(908, 45462, 0.13782, 0.13783, 0.13775, 0.13775, 0),
    @(909, 45463, 0.13775, 0.13775, 0.13769, 0.13769, 0),
    @(910, 45464, 0.1377, 0.13771, 0.13767, 0.13767, 0),
    @(911, 45467, 0.13768, 0.13776, 0.13768, 0.13773, 0),
    @(912, 45468, 0.13772, 0.13775, 0.13765, 0.13765, 0),
    @(913, 45469, 0.13765, 0.13766, 0.13761, 0.13761, 0),
    @(914, 45470, 0.13761, 0.13762, 0.13753, 0.13753, 0),
    @(915, 45471, 0.13755, 0.13766, 0.13755, 0.13756, 0),
    @(916, 45474, 0.13756, 0.13757, 0.13754, 0.13754, 0),
    @(917, 45475, 0.13754, 0.13755, 0.13749, 0.13749, 0),
    @(918, 45476, 0.13749, 0.13755, 0.13749, 0.13751, 0),
    @(919, 45477, 0.13751, 0.13757, 0.13751, 0.13756, 0),
    @(920, 45478, 0.13757, 0.13761, 0.13753, 0.13753, 0),
    @(921, 45481, 0.13755, 0.13756, 0.13754, 0.13754, 0),
    @(922, 45482, 0.13754, 0.13755, 0.13746, 0.13746, 0),
    @(923, 45483, 0.13746, 0.13748, 0.1374, 0.1374, 0),
    @(924, 45484, 0.13741, 0.13788, 0.13741, 0.13774, 0),
    @(925, 45485, 0.13774, 0.13791, 0.13769, 0.13787, 0),
    @(926, 45488, 0.13789, 0.13789, 0.13771, 0.13772, 0),
    @(927, 45489, 0.13772, 0.13837, 0.13754, 0.13754, 0),
    @(928, 45490, 0.13754, 0.13885, 0.13754, 0.13769, 0),
    @(929, 45491, 0.13769, 0.13781, 0.13769, 0.13769, 0),
    @(930, 45492, 0.1377, 0.13772, 0.13751, 0.13751, 0),
    @(931, 45495, 0.13752, 0.13755, 0.13744, 0.13744, 0),
    @(932, 45496, 0.13744, 0.13746, 0.13742, 0.13742, 0),
    @(933, 45497, 0.13742, 0.13766, 0.13742, 0.13765, 0),
    @(934, 45498, 0.13765, 0.13865, 0.13765, 0.13825, 0),
    @(935, 45499, 0.13825, 0.13825, 0.13789, 0.13789, 0),
    @(936, 45502, 0.13789, 0.1379, 0.13771, 0.13771, 0),
    @(937, 45503, 0.13771, 0.13792, 0.13769, 0.13787, 0),
    @(938, 45504, 0.13787, 0.13855, 0.13787, 0.13848, 0),
    @(939, 45505, 0.13846, 0.13864, 0.138, 0.138, 0),
    @(940, 45506, 0.13802, 0.1397, 0.13802, 0.13962, 0),
    @(941, 45509, 0.13962, 0.14135, 0.13962, 0.14021, 0),
    @(942, 45510, 0.14021, 0.1406, 0.1397, 0.1397, 0),
    @(943, 45511, 0.1397, 0.14024, 0.13919, 0.13931, 0),
    @(944, 45512, 0.13931, 0.13969, 0.13931, 0.13934, 0),
    @(945, 45513, 0.13934, 0.13963, 0.13934, 0.13951, 0),
    @(946, 45516, 0.1395, 0.13994, 0.13927, 0.13932, 0),
    @(947, 45517, 0.13934, 0.14, 0.13927, 0.13972, 0),
    @(948, 45518, 0.14048, 0.14051, 0.13983, 0.14008, 0),
    @(949, 45519, 0.14007, 0.14035, 0.13964, 0.13978, 0)
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]

    # Match the existing date-column formatting (style index used by column A)
    $styleSource.Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

Write-Output ("New dimension: " + $ws.UsedRange.Address())
